# Adds an explicit <w:pageBreakBefore w:val="0"/> to every paragraph's
# paragraph-properties (w:pPr) in the document body, and to every paragraph
# style definition that already carries paragraph formatting (the Heading
# 1-6, Title and Subtitle styles). This mirrors a Google Docs -> OOXML
# export quirk where the (false) "page break before" setting is written out
# explicitly instead of being omitted.

$d = $word.ActiveDocument

# 1. Every paragraph in the main document body.
foreach ($p in $d.Paragraphs) {
    $p.Format.PageBreakBefore = 0
}

# 2. The paragraph styles that define paragraph formatting (keepNext /
#    keepLines / spacing) - Heading 1-6, Title, Subtitle.
$styleNames = @("Heading 1", "Heading 2", "Heading 3", "Heading 4", "Heading 5", "Heading 6", "Title", "Subtitle")
foreach ($name in $styleNames) {
    $s = $d.Styles.Item($name)
    $s.ParagraphFormat.PageBreakBefore = 0
}
